# feat: add 2022-Q1 data
#
# The workbook has two sheets: "2020-Q4" (fund holdings detail) and
# "总计" (summary of holdings per quarter). This change:
#   1. Turns the current "总计" sheet into the new "2022-Q1" detail sheet
#      (same column layout as "2020-Q4") with its own fund holdings data.
#   2. Adds a brand-new "总计" sheet (placed after "2022-Q1") that keeps the
#      previous summary row (2020-Q4) and adds a new row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(2)   # currently named "总计"

# Duplicate the summary sheet first (right after itself) while it still has
# its original content/formatting - this duplicate will become the new
# "总计" sheet. The original will be turned into "2022-Q1".
$summary.Copy($null, $summary)

$quarter = $wb.Worksheets.Item(2)     # the original "总计" sheet
$newSummary = $wb.Worksheets.Item(3)  # the freshly created duplicate

$quarter.Name = "2022-Q1"
$newSummary.Name = "总计"

# ---------------------------------------------------------------------
# 1) Rebuild the "2022-Q1" sheet with the fund holdings layout/content
# ---------------------------------------------------------------------
$quarter.Range("A1:D2").ClearContents()

# Extend the header style (currently only on B1:D1) across to H1 by
# copying the formatting (not the value) of an already-styled header cell.
$quarter.Range("B1").Copy()
$quarter.Range("E1:H1").PasteSpecial(-4122)

$quarter.Range("B1").Value = "基金代码"
$quarter.Range("C1").Value = "基金名称"
$quarter.Range("D1").Value = "基金规模"
$quarter.Range("E1").Value = "股票总仓位"
$quarter.Range("F1").Value = "仓位占比"
$quarter.Range("G1").Value = "持有市值(亿元)"
$quarter.Range("H1").Value = "仓位排名"

$quarter.Range("A2").Value = 0

$quarter.Range("B2").NumberFormat = "@"
$quarter.Range("B2").Value = "004250"
$quarter.Range("B2").ClearFormats()

$quarter.Range("C2").Value = "银河量化优选混合"

$quarter.Range("D2").NumberFormat = "@"
$quarter.Range("D2").Value = "0.39"
$quarter.Range("D2").ClearFormats()

$quarter.Range("E2").NumberFormat = "@"
$quarter.Range("E2").Value = "80.03"
$quarter.Range("E2").ClearFormats()

$quarter.Range("F2").NumberFormat = "@"
$quarter.Range("F2").Value = "2.23"
$quarter.Range("F2").ClearFormats()

$quarter.Range("G2").NumberFormat = "@"
$quarter.Range("G2").Value = "0.0087"
$quarter.Range("G2").ClearFormats()

$quarter.Range("H2").Value = 2

# ---------------------------------------------------------------------
# 2) Update the new "总计" sheet: insert a 2022-Q1 row above the existing
#    2020-Q4 row (which shifts down to row 3 while keeping its values).
# ---------------------------------------------------------------------
$newSummary.Rows(2).Insert()
$newSummary.Range("B2:D2").ClearFormats()

$newSummary.Range("A3").Copy()
$newSummary.Range("A2").PasteSpecial(-4122)

$newSummary.Range("A2").Value = 0
$newSummary.Range("B2").Value = "2022-Q1"
$newSummary.Range("C2").Value = 1
$newSummary.Range("D2").Value = 0.01

# The row that got shifted down (2020-Q4) now has index 1 (0-based row
# index kept in column A), not its previous 0.
$newSummary.Range("A3").Value = 1

# Restore "2020-Q4" as the active/selected sheet, matching the original
# workbook state (it was the only sheet with tabSelected="1").
$wb.Worksheets.Item(1).Activate()
